$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.95500663279169
$ws.Range("C2").Value = 10.27372646331787
$ws.Range("D2").Value = 6.62
$ws.Range("E2").Value = 0.8971493667929927

$ws.Range("B3").Value = 14.76877601546984
$ws.Range("C3").Value = 11.66212177276611
$ws.Range("D3").Value = 7.37
$ws.Range("E3").Value = 0.9784173339000368

$ws.Range("B4").Value = 2.726994696085647
$ws.Range("C4").Value = 2.046205520629883
$ws.Range("D4").Value = 1.33
$ws.Range("E4").Value = 0.9505758440998126

$ws.Range("B5").Value = 2.726892314736134
$ws.Range("C5").Value = 2.046073913574219
$ws.Range("D5").Value = 1.33
$ws.Range("E5").Value = 0.9505798530221522

$ws.Range("B6").Value = 12.58162679321207
$ws.Range("C6").Value = 7.802098751068115
$ws.Range("D6").Value = 7.08
$ws.Range("E6").Value = 0.857517469373574

$ws.Range("B7").Value = 8.572025541679576
$ws.Range("C7").Value = 6.120454788208008
$ws.Range("D7").Value = 3.33
$ws.Range("E7").Value = 0.9302336006405363

$ws.Range("B8").Value = 9.165262519335704
$ws.Range("C8").Value = 6.709384441375732
$ws.Range("D8").Value = 5.1
$ws.Range("E8").Value = 0.9331783782381344

$ws.Range("B9").Value = 13.90684250844481
$ws.Range("C9").Value = 9.659384727478027
$ws.Range("D9").Value = 9.970000000000001
$ws.Range("E9").Value = 0.8511621495469708

$ws.Range("B10").Value = 60.14826448621045
$ws.Range("C10").Value = 33.1715202331543
$ws.Range("D10").Value = 28.94
$ws.Range("E10").Value = -0.5487403543442707

$wb.Save()
